{"js": "// Replace each \"AxB=\" arithmetic-drill expression in the document with its\n// new value. Every source expression is unique in the document, so a plain\n// text search-and-replace (no wildcards needed) is unambiguous.\nconst replacements = [\n  [\"741\u00d78=\", \"605\u00d74=\"],\n  [\"932\u00d73=\", \"573\u00d79=\"],\n  [\"499\u00d73=\", \"241\u00d79=\"],\n  [\"226\u00d76=\", \"729\u00d75=\"],\n  [\"576\u00d77=\", \"953\u00d77=\"],\n  [\"495\u00d76=\", \"364\u00d77=\"],\n  [\"769\u00d77=\", \"217\u00d78=\"],\n  [\"571\u00d78=\", \"578\u00d72=\"],\n  [\"197\u00d73=\", \"190\u00d76=\"],\n  [\"395\u00d72=\", \"952\u00d79=\"],\n  [\"407\u00d75=\", \"460\u00d73=\"],\n  [\"233\u00d76=\", \"786\u00d77=\"],\n  [\"992\u00d72=\", \"163\u00d78=\"],\n  [\"309\u00d74=\", \"556\u00d74=\"],\n  [\"336\u00d76=\", \"214\u00d74=\"],\n  [\"962\u00d78=\", \"353\u00d76=\"],\n  [\"963\u00d72=\", \"769\u00d76=\"],\n  [\"189\u00d75=\", \"365\u00d78=\"],\n  [\"191\u00d73=\", \"309\u00d76=\"],\n  [\"640\u00d74=\", \"660\u00d75=\"],\n  [\"109\u00d76=\", \"890\u00d76=\"],\n  [\"571\u00d73=\", \"554\u00d74=\"],\n  [\"920\u00d78=\", \"379\u00d77=\"],\n  [\"784\u00d76=\", \"498\u00d72=\"],\n  [\"190\u00d79=\", \"112\u00d73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each \"AxB=\" arithmetic-drill expression in the document with its\n# new value. Every source expression is unique in the document, so a plain\n# Find/Replace (wdReplaceAll, no wildcards) is unambiguous for each pair.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"741\u00d78=\", \"605\u00d74=\"),\n    @(\"932\u00d73=\", \"573\u00d79=\"),\n    @(\"499\u00d73=\", \"241\u00d79=\"),\n    @(\"226\u00d76=\", \"729\u00d75=\"),\n    @(\"576\u00d77=\", \"953\u00d77=\"),\n    @(\"495\u00d76=\", \"364\u00d77=\"),\n    @(\"769\u00d77=\", \"217\u00d78=\"),\n    @(\"571\u00d78=\", \"578\u00d72=\"),\n    @(\"197\u00d73=\", \"190\u00d76=\"),\n    @(\"395\u00d72=\", \"952\u00d79=\"),\n    @(\"407\u00d75=\", \"460\u00d73=\"),\n    @(\"233\u00d76=\", \"786\u00d77=\"),\n    @(\"992\u00d72=\", \"163\u00d78=\"),\n    @(\"309\u00d74=\", \"556\u00d74=\"),\n    @(\"336\u00d76=\", \"214\u00d74=\"),\n    @(\"962\u00d78=\", \"353\u00d76=\"),\n    @(\"963\u00d72=\", \"769\u00d76=\"),\n    @(\"189\u00d75=\", \"365\u00d78=\"),\n    @(\"191\u00d73=\", \"309\u00d76=\"),\n    @(\"640\u00d74=\", \"660\u00d75=\"),\n    @(\"109\u00d76=\", \"890\u00d76=\"),\n    @(\"571\u00d73=\", \"554\u00d74=\"),\n    @(\"920\u00d78=\", \"379\u00d77=\"),\n    @(\"784\u00d76=\", \"498\u00d72=\"),\n    @(\"190\u00d79=\", \"112\u00d73=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
